$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 12.07002266666666
$ws.Range("H2").Value = 36.21006799999999
$ws.Range("I2").Value = 0.7601982364861632
$ws.Range("J2").Value = 0.7601982364861634
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.464265666666667
$ws.Range("N2").Value = 4.392797
$ws.Range("O2").Value = 0.02620474750556022
$ws.Range("P2").Value = 0.02620474750556022
$ws.Range("Q2").Value = 17.67371978668844
$ws.Range("R2").Value = 159.0634780801959
$ws.Range("S2").Value = 0.01992080284129207
$ws.Range("T2").Value = 0.01992080284129207

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 12.07002266666666
$ws.Range("H3").Value = 36.21006799999999
$ws.Range("I3").Value = 0.7601982364861632
$ws.Range("J3").Value = 0.7601982364861634
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8090393333333333
$ws.Range("N3").Value = 2.427118
$ws.Range("O3").Value = 0.01447870556190061
$ws.Range("P3").Value = 0.01447870556190061
$ws.Range("Q3").Value = 9.765123091558221
$ws.Range("R3").Value = 87.88610782402398
$ws.Range("S3").Value = 0.01100668643475925
$ws.Range("T3").Value = 0.01100668643475925

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.07002266666666
$ws.Range("H4").Value = 36.21006799999999
$ws.Range("I4").Value = 0.7601982364861632
$ws.Range("J4").Value = 0.7601982364861634
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 53.60457233333333
$ws.Range("N4").Value = 160.813717
$ws.Range("O4").Value = 0.9593165469325392
$ws.Range("P4").Value = 0.9593165469325391
$ws.Range("Q4").Value = 647.008403100306
$ws.Range("R4").Value = 5823.075627902755
$ws.Range("S4").Value = 0.729270747210112
$ws.Range("T4").Value = 0.729270747210112

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.308268
$ws.Range("H5").Value = 3.924804
$ws.Range("I5").Value = 0.08239777620284613
$ws.Range("J5").Value = 0.08239777620284613
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.464265666666667
$ws.Range("N5").Value = 4.392797
$ws.Range("O5").Value = 0.02620474750556022
$ws.Range("P5").Value = 0.02620474750556022
$ws.Range("Q5").Value = 1.915651915198667
$ws.Range("R5").Value = 17.240867236788
$ws.Range("S5").Value = 0.002159212920415242
$ws.Range("T5").Value = 0.002159212920415242

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.308268
$ws.Range("H6").Value = 3.924804
$ws.Range("I6").Value = 0.08239777620284613
$ws.Range("J6").Value = 0.08239777620284613
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8090393333333333
$ws.Range("N6").Value = 2.427118
$ws.Range("O6").Value = 0.01447870556190061
$ws.Range("P6").Value = 0.01447870556190061
$ws.Range("Q6").Value = 1.058440270541333
$ws.Range("R6").Value = 9.525962434872
$ws.Range("S6").Value = 0.00119301314059639
$ws.Range("T6").Value = 0.00119301314059639

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.308268
$ws.Range("H7").Value = 3.924804
$ws.Range("I7").Value = 0.08239777620284613
$ws.Range("J7").Value = 0.08239777620284613
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 53.60457233333333
$ws.Range("N7").Value = 160.813717
$ws.Range("O7").Value = 0.9593165469325392
$ws.Range("P7").Value = 0.9593165469325391
$ws.Range("Q7").Value = 70.12914663738533
$ws.Range("R7").Value = 631.162319736468
$ws.Range("S7").Value = 0.07904555014183451
$ws.Range("T7").Value = 0.0790455501418345

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.499176666666667
$ws.Range("H8").Value = 7.49753
$ws.Range("I8").Value = 0.1574039873109905
$ws.Range("J8").Value = 0.1574039873109906
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.464265666666667
$ws.Range("N8").Value = 4.392797
$ws.Range("O8").Value = 0.02620474750556022
$ws.Range("P8").Value = 0.02620474750556022
$ws.Range("Q8").Value = 3.659458587934445
$ws.Range("R8").Value = 32.93512729141
$ws.Range("S8").Value = 0.004124731743852912
$ws.Range("T8").Value = 0.004124731743852912

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.499176666666667
$ws.Range("H9").Value = 7.49753
$ws.Range("I9").Value = 0.1574039873109905
$ws.Range("J9").Value = 0.1574039873109906
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8090393333333333
$ws.Range("N9").Value = 2.427118
$ws.Range("O9").Value = 0.01447870556190061
$ws.Range("P9").Value = 0.01447870556190061
$ws.Range("Q9").Value = 2.021932224282222
$ws.Range("R9").Value = 18.19739001854
$ws.Range("S9").Value = 0.002279005986544972
$ws.Range("T9").Value = 0.002279005986544972

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.499176666666667
$ws.Range("H10").Value = 7.49753
$ws.Range("I10").Value = 0.1574039873109905
$ws.Range("J10").Value = 0.1574039873109906
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 53.60457233333333
$ws.Range("N10").Value = 160.813717
$ws.Range("O10").Value = 0.9593165469325392
$ws.Range("P10").Value = 0.9593165469325391
$ws.Range("Q10").Value = 133.9672964021122
$ws.Range("R10").Value = 1205.70566761901
$ws.Range("S10").Value = 0.1510002495805927
$ws.Range("T10").Value = 0.1510002495805927

Write-Output "Done"
